# Update remaining template's ontology terms
#
# 1. On the main annotation sheet (4COM04_GenomeAssembly), fill in a new
#    data row describing the "Next generation sequencing instrument model"
#    parameter with the term "Illumina Genome Analyzer IIx" and its
#    ontology reference (NFDI4PSO).
# 2. Bump the template version on the SwateTemplateMetadata sheet from
#    1.1.4 to 1.1.5, and make that sheet the active one.

$wb = $excel.ActiveWorkbook

# --- 1. Main annotation sheet: add row 2 values -----------------------
$ws1 = $wb.Worksheets.Item("4COM04_GenomeAssembly")
$ws1.Range("N2").Value = "Illumina Genome Analyzer IIx"
$ws1.Range("O2").Value = "NFDI4PSO"
$ws1.Range("P2").Value = "http://purl.obolibrary.org/obo/NFDI4PSO_1000038"

# Widen the (hidden) ontology-accession column so the new, longer URL
# fits, matching the sheet's other auto-fitted accession columns.
$ws1.Columns.Item(16).ColumnWidth = 45.666666666666664

# --- 2. Metadata sheet: bump version string ----------------------------
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")
# Leading apostrophe keeps the cell's existing "stored as text"
# (quote-prefixed) formatting instead of resetting it.
$ws2.Range("B3").Value = "'1.1.5"

# SwateTemplateMetadata becomes the active/visible sheet.
$ws2.Activate()
